$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the bottom of the tracker table (row 90), shifting
# down and inheriting the formatting of the row above (mirrors what Excel
# does when a user inserts a row from the UI).
$ws.Range("A90:E90").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Fill in the new entry: 102. Binary Tree Level Order Traversal
$ws.Range("A90").Value = 102
$ws.Range("B90").Value = "Binary Tree Level Order Traversal"
$ws.Range("C90").Value = "Medium"
$ws.Range("D90").Value = "BFS,level order traversal"
$ws.Range("E90").Value = 45785

# Restore the saved selection state from the workbook view
$ws.Range("B92").Select()
